$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of robustness-of-bias (RoB) data extracted and merged into the
# studylist, appended below the existing data (rows 131-140).
# Columns used: A, C, D, E, F, G, H, I, J, K, L, M, O, P, Q, R, S
# (columns B "Review" and N "Study" are intentionally left blank, matching
# the pattern of the existing "Rob_HPP_J" rows.)

$newRows = @(
    @("Rob_HPP_J","Barrowclough(2015)","longitudinal prospective","1","1","1","0","2","1","1","7","0","Johanna","Barrowclough","Barrowclough","2015","barrowclough_2015"),
    @("Rob_HPP_J","Barrowclough(2013)","longitudinal prospective","1","1","1","0","2","1","1","7","0","Johanna","Barrowclough","Barrowclough","2013","barrowclough_2013"),
    @("Rob_HPP_J","Fond(2019)","longitudinal prospective","1","1","1","0","2","0","1","6","0","Johanna","Fond","Fond","2019","fond_2019"),
    @("Rob_HPP_J","Foti(2010)","longitudinal prospective","1","1","1","0","2","1","1","8","1","Johanna","Foti","Foti","2010","foti_2010"),
    @("Rob_HPP_J","Baeza(2009)","longitudinal prospective","1","1","0","1","2","1","1","8","1","Johanna","Baeza","Baeza","2009","baeza_2009"),
    @("Rob_HPP_J","Buchy(2015)","longitudinal prospective","1","1","1","1","1","1","1","7","0","Johanna","Buchy","Buchy","2015","buchy_2015"),
    @("Rob_HPP_J","Zammit(2011)","longitudinal prospective","1","1","0","1","2","1","1","8","1","Johanna","Zammit","Zammit","2011","zammit_2011"),
    @("Rob_HPP_J","Arseneault(2002)","longitudinal prospective","1","1","0","0","2","1","1","7","1","Johanna","Arseneault","Arseneault","2002","arseneault_2002"),
    @("Rob_HPP_J","Bechtold(2016)","longitudinal prospective","1","1","0","0","2","0","1","6","1","Johanna","Bechtold","Bechtold","2016","bechtold_2016"),
    @("Rob_HPP_J","Dragt(2011)","longitudinal prospective","1","1","1","1","2","1","1","9","1","Johanna","Dragt","Dragt","2011","dragt_2011")
)

# Columns in the order that values appear within each row above.
$colLetters = @("A","C","D","E","F","G","H","I","J","K","L","M","O","P","Q","R","S")

$startRow = 131
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "$($colLetters[$j])$rowNum"
        $cell = $ws.Range($cellRef)
        # Force text storage so numeric-looking values (e.g. "1", "2015")
        # are written as strings, matching the rest of the sheet where
        # every value is an inline/shared string rather than a number.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$j]
    }
}
